$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "Förändrad" (column C) date bumped from 2023-09-15 (45184) to 2023-09-17 (45186)
#    for every data row (2 through 72).
$ws.Range("C2:C72").Value = 45186

# 2) The HYPERLINK() formulas in columns S, T, V, W, X, Y (rows 2-16, the rows that
#    have those link columns populated) gained a second "friendly name" argument
#    equal to the report's "Beteckning" (column A) value.
$folder = @{ "S" = "artfynd"; "T" = "kartor"; "V" = "klagomål"; "W" = "klagomålsmail"; "X" = "tillsyn"; "Y" = "tillsynsmail" }
$ext    = @{ "S" = "xlsx"; "T" = "png"; "V" = "docx"; "W" = "docx"; "X" = "docx"; "Y" = "docx" }

for ($row = 2; $row -le 16; $row++) {
    $beteckning = $ws.Cells.Item($row, 1).Value2
    foreach ($col in @("S", "T", "V", "W", "X", "Y")) {
        $url = "https://klasma.github.io/Logging_VARMDO/" + $folder[$col] + "/" + $beteckning + "." + $ext[$col]
        $ws.Range("$col$row").Formula = "=HYPERLINK(`"$url`", `"$beteckning`")"
    }
}
